# Add a "Comments" column (column E) to the four history sheets, matching
# the header style already used in row 1 of each sheet, then leave the
# "Withdraw History" sheet active with E1 selected.

$wb = $excel.ActiveWorkbook

$sheetNames = @("Withdraw History", "Deposit History", "Transfer History", "Absolute History")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("E1").Value = "Comments"
    $ws.Range("E1").Select()
}

# Final active sheet/selection: "Withdraw History", cell E1.
$wb.Worksheets.Item("Withdraw History").Activate()
$wb.Worksheets.Item("Withdraw History").Range("E1").Select()
